# Wilke spinal disc pressure validation study — update measured AMS forces
# (column C) after fixing a transversus symmetry bug. The adjacent ratio
# column D (=Cn/C5) and the chart series that plots it recompute on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("C2").Value = 99.737080000000006
$ws.Range("C3").Value = 282.56619999999998
$ws.Range("C4").Value = 570.31399999999996
$ws.Range("C5").Value = 663.96690000000001
$ws.Range("C6").Value = 1392.172
$ws.Range("C7").Value = 1388.2860000000001
$ws.Range("C8").Value = 2892.1970000000001
$ws.Range("C9").Value = 2474.6509999999998

# Leave the selection where the author's last interaction left it.
$ws.Range("G14").Select()
